$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- New row 9: Form_URL / hyperlink to JotForm ---
$ws.Range("B9").Value = "https://form.jotform.com/250753158727665"
$ws.Range("A9").Value = "Form_URL"
$ws.Hyperlinks.Add($ws.Range("B9"), "https://form.jotform.com/250753158727665")
$ws.Range("B9").Font.Underline = 0

# --- New row 10: Colaborador / Gabriel Ballone ---
$ws.Range("A10").Value = "Colaborador"
$ws.Range("B10").Value = "Gabriel Ballone"

# --- New (blank, formatted) cell at B12, matching the style already used at C10 ---
$ws.Range("C10").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection / view state ---
$ws.Range("B12").Select()
